$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying the style used by the existing
# header cells (e.g. H1) so they match formatting-wise.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Add the new data values in row 2
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
